# "hydrogen and footer fix"
#
# 1. The stray "_GoBack" bookmark that used to sit right after the
#    "электро-" run (before "водородных автомобилей...") is removed.
# 2. A new "_GoBack" bookmark is added spanning from the start of
#    "Сделать топливо для машины..." through the end of
#    "...к 2077 году." (i.e. wrapping the whole last "paragraph" of
#    body text, but stopping before the trailing line break run).

$d = $word.ActiveDocument

# --- 1. Remove the old bookmark --------------------------------------------
# (Bookmarks("_GoBack").Delete() is a safe no-op if it is somehow missing.)
$d.Bookmarks("_GoBack").Delete()

# --- 2. Work out the new bookmark's start/end character offsets -----------
$startRange = $d.Content
$startRange.Find.Execute("Сделать топливо для машины из воды") | Out-Null
$startPos = $startRange.Start

$endRange = $d.Content
$endRange.Find.Execute("будут стоять на одной ценовой ступеньке с ДВС, например, к 2077 году.") | Out-Null
$endPos = $endRange.End

# --- 3. Re-create "_GoBack" over that span ---------------------------------
$bookmarkRange = $d.Range($startPos, $endPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
